# "Generate Report for Archive"
# - Status changes from "Ready for handoff" to "In Translation" on every
#   sheet that surfaces it (Overview!E2:F2, zh-cn!C2, de-de!C2).
# - The zh-cn/de-de status column (and the matching columns on Overview)
#   are narrowed to the new auto-fit width.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Update the status text wherever "Ready for handoff" is shown.
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# Narrow the status columns to match the new content width.
# (Columns.Item takes a 1-based numeric column index here.)
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5

Write-Output "Report regenerated for archive."
